$wb = $excel.ActiveWorkbook

# Sheet 1: updates to column F ("想去人数")
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 119
$ws.Cells.Item(4, 6).Value = 880
$ws.Cells.Item(5, 6).Value = 1057
$ws.Cells.Item(6, 6).Value = 1542
$ws.Cells.Item(8, 6).Value = 652
$ws.Cells.Item(9, 6).Value = 11907
$ws.Cells.Item(11, 6).Value = 2137
$ws.Cells.Item(13, 6).Value = 246
$ws.Cells.Item(16, 6).Value = 1217
$ws.Cells.Item(17, 6).Value = 185
$ws.Cells.Item(18, 6).Value = 260
$ws.Cells.Item(19, 6).Value = 751
$ws.Cells.Item(20, 6).Value = 666
$ws.Cells.Item(21, 6).Value = 285
$ws.Cells.Item(22, 6).Value = 2912
$ws.Cells.Item(23, 6).Value = 746
$ws.Cells.Item(24, 6).Value = 3802
$ws.Cells.Item(25, 6).Value = 3802
$ws.Cells.Item(26, 6).Value = 1084
$ws.Cells.Item(31, 6).Value = 1017
$ws.Cells.Item(33, 6).Value = 89
$ws.Cells.Item(37, 6).Value = 22
$ws.Cells.Item(38, 6).Value = 4352
$ws.Cells.Item(39, 6).Value = 12
$ws.Cells.Item(40, 6).Value = 4491
$ws.Cells.Item(41, 6).Value = 5512
$ws.Cells.Item(44, 6).Value = 54
$ws.Cells.Item(45, 6).Value = 168
$ws.Cells.Item(46, 6).Value = 286
$ws.Cells.Item(47, 6).Value = 74
$ws.Cells.Item(48, 6).Value = 38
$ws.Cells.Item(49, 6).Value = 4102
$ws.Cells.Item(50, 6).Value = 115

# Sheet 2: updates to column F ("想去人数")
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(3, 6).Value = 4168
$ws.Cells.Item(10, 6).Value = 108
$ws.Cells.Item(12, 6).Value = 845

# Sheet 3: updates to column F ("想去人数")
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 6).Value = 431

# Sheet 4: updates to column F ("想去人数")
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(3, 6).Value = 431
$ws.Cells.Item(6, 6).Value = 880
$ws.Cells.Item(7, 6).Value = 1057
$ws.Cells.Item(8, 6).Value = 1542
$ws.Cells.Item(10, 6).Value = 652
$ws.Cells.Item(11, 6).Value = 11907
$ws.Cells.Item(12, 6).Value = 2137
$ws.Cells.Item(14, 6).Value = 246
$ws.Cells.Item(15, 6).Value = 1217
$ws.Cells.Item(16, 6).Value = 185
$ws.Cells.Item(17, 6).Value = 260
$ws.Cells.Item(18, 6).Value = 4168
$ws.Cells.Item(19, 6).Value = 751
$ws.Cells.Item(20, 6).Value = 285
$ws.Cells.Item(21, 6).Value = 746
$ws.Cells.Item(22, 6).Value = 3802
$ws.Cells.Item(23, 6).Value = 1084
$ws.Cells.Item(29, 6).Value = 1017
$ws.Cells.Item(31, 6).Value = 89
$ws.Cells.Item(34, 6).Value = 22
$ws.Cells.Item(35, 6).Value = 4491
$ws.Cells.Item(38, 6).Value = 168
$ws.Cells.Item(39, 6).Value = 286
$ws.Cells.Item(43, 6).Value = 74
$ws.Cells.Item(44, 6).Value = 38
$ws.Cells.Item(45, 6).Value = 4102
$ws.Cells.Item(50, 6).Value = 115
